# Applies the heart_0531 edit: flips (toggles 0<->1) the label column N
# for rows 408 through 938 (inclusive), and updates the sheet's active
# selection to O408 (matching the scrolled/selected view in the edited file).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 408
$lastRow  = 938
$col      = 14   # column N

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $old  = $cell.Value2
    $new  = 1 - $old
    $cell.Value = $new
}

# Reflect the saved view state: window scrolled so row 401 is at the top,
# with O408 as the active/selected cell.
$ws.Range("O408").Select()
$excel.ActiveWindow.ScrollRow = 401
